# Updated cryptos list with GitHub Actions
#
# Note: Price (D) and Volume(1h) (E) columns are plain text in this sheet
# (not numbers), e.g. "23.455.22", "0.9979", "  -0.94%  ". Several of these
# look like numbers/dates to Excel's auto-detection, so every literal is
# prefixed with a leading apostrophe ('') - PowerShell's escaped single
# quote inside a single-quoted string - which Excel treats as a text
# qualifier: it forces the cell to stay text without altering NumberFormat.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''23.455.22'
$ws.Range('E2').Value = '''  -0.94%  '
$ws.Range('D3').Value = '''1.648.36'
$ws.Range('E3').Value = '''  -0.38%  '
$ws.Range('D4').Value = '''0.9979'
$ws.Range('E4').Value = '''  -0.43%  '
$ws.Range('D5').Value = '''0.9990'
$ws.Range('E5').Value = '''  -0.28%  '
$ws.Range('D6').Value = '''299.73'
$ws.Range('E6').Value = '''  -1.10%  '
$ws.Range('D7').Value = '''0.3800'
$ws.Range('E7').Value = '''  -0.93%  '
$ws.Range('D8').Value = '''50.42'
$ws.Range('E8').Value = '''  -1.30%  '
$ws.Range('D9').Value = '''0.3495'
$ws.Range('E9').Value = '''  -3.08%  '
$ws.Range('D10').Value = '''0.08075'
$ws.Range('D11').Value = '''1.219'
$ws.Range('D12').Value = '''0.9979'
$ws.Range('E12').Value = '''  -0.42%  '
$ws.Range('D13').Value = '''22.07'
$ws.Range('E13').Value = '''  -1.49%  '
$ws.Range('D14').Value = '''6.324'
$ws.Range('E14').Value = '''  -2.03%  '
$ws.Range('D15').Value = '''7.285'
$ws.Range('E15').Value = '''  -2.10%  '
$ws.Range('D16').Value = '''0.00001216'
$ws.Range('E16').Value = '''  -0.46%  '
$ws.Range('D17').Value = '''1.647.99'
$ws.Range('E17').Value = '''  -0.40%  '
$ws.Range('D18').Value = '''94.90'
$ws.Range('E18').Value = '''  -2.78%  '
$ws.Range('D19').Value = '''0.06970'
$ws.Range('E19').Value = '''  -0.94%  '
$ws.Range('D20').Value = '''6.638'
$ws.Range('E20').Value = '''  -1.97%  '
$ws.Range('D21').Value = '''17.39'
$ws.Range('E21').Value = '''  -0.92%  '
$ws.Range('D22').Value = '''0.9983'
$ws.Range('E22').Value = '''  -0.31%  '
$ws.Range('E23').Value = '''  -2.03%  '
$ws.Range('D24').Value = '''23.453.25'
$ws.Range('E24').Value = '''  -1.00%  '
$ws.Range('D25').Value = '''2.436'
$ws.Range('E25').Value = '''  -1.96%  '
$ws.Range('D26').Value = '''2.973'
$ws.Range('E26').Value = '''  -1.66%  '
$ws.Range('D27').Value = '''21.05'
$ws.Range('E27').Value = '''  -0.94%  '
$ws.Range('D28').Value = '''149.76'
$ws.Range('E28').Value = '''  -2.28%  '
$ws.Range('D29').Value = '''5.180'
$ws.Range('E29').Value = '''  -1.07%  '
$ws.Range('D30').Value = '''131.76'
$ws.Range('E30').Value = '''  -1.86%  '
$ws.Range('D31').Value = '''1.825.70'
$ws.Range('E31').Value = '''  -0.76%  '
$ws.Range('D32').Value = '''6.849'
$ws.Range('E32').Value = '''  -3.38%  '
$ws.Range('D33').Value = '''2.131'
$ws.Range('E33').Value = '''  -5.51%  '
$ws.Range('E34').Value = '''  -6.75%  '
$ws.Range('D35').Value = '''0.9905'
$ws.Range('E35').Value = '''  -6.42%  '
$ws.Range('D36').Value = '''0.02687'
$ws.Range('E36').Value = '''  -4.22%  '
$ws.Range('D37').Value = '''0.08791'
$ws.Range('E37').Value = '''  -0.17%  '
$ws.Range('D38').Value = '''0.2424'
$ws.Range('E38').Value = '''  -3.29%  '
$ws.Range('D39').Value = '''5.891'
$ws.Range('E39').Value = '''  -3.03%  '
$ws.Range('D40').Value = '''0.06846'
$ws.Range('E40').Value = '''  -1.99%  '
$ws.Range('D41').Value = '''12.77'
$ws.Range('E41').Value = '''  -1.90%  '
$ws.Range('D42').Value = '''0.6840'
$ws.Range('E42').Value = '''  -2.07%  '
$ws.Range('D43').Value = '''1.287'
$ws.Range('E43').Value = '''  -3.76%  '
$ws.Range('D44').Value = '''15.53'
$ws.Range('E44').Value = '''  -2.47%  '
$ws.Range('D45').Value = '''0.9979'
$ws.Range('E45').Value = '''  -0.35%  '
$ws.Range('D46').Value = '''0.6359'
$ws.Range('E46').Value = '''  -2.18%  '
$ws.Range('D47').Value = '''2.243'
$ws.Range('E47').Value = '''  -2.49%  '
$ws.Range('D48').Value = '''3.910'
$ws.Range('E48').Value = '''  -1.39%  '
$ws.Range('D51').Value = '''1.224'
$ws.Range('E51').Value = '''  +2.60%  '

# Row 49/50: Quant now ranks above Cronos, with updated price/volume
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '''127.32'
$ws.Range('E49').Value = '''  -0.63%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.07691'
$ws.Range('E50').Value = '''  -2.50%  '
